# Update "想去人数" (number of people interested) counts for two events
# that appear on both the "展览" sheet and the aggregated "全部类型" sheet.
#   F2: 302 -> 304
#   F4: 1221 -> 1226

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 304
    $ws.Range("F4").Value = 1226
}
